# LimitManagement.xlsx update
# Adds two new limit rows (Funds Transfer to HBL Account / Funds Transfer to
# Other Banks Account) and repurposes the existing row 2 as the "Utility
# Bills and Other Payments" limit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows directly below the existing data row (row 2), so the
# sheet grows from rows 1:2 to rows 1:4.
$ws.Rows("3:4").Insert()

# Populate the columns that are identical across every data row (A, D-I) for
# the two new rows, reusing the exact same text already used on row 2.
$ws.Range("A3").Value = "As A user I want to verify Limit Reduction"
$ws.Range("D3").Value = "Non Financial"
$ws.Range("E3").Value = "SELECT AC.DESCRIPTION FROM DC_TRANSACTION_ACTIVITY_CONFIG AC WHERE AC.TRANSACTION_TYPE_ID =(SELECT DT.TRANSACTION_TYPE_ID FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("F3").Value = "SELECT DT.IVR_ATTRIBUTE2 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("G3").Value = "SELECT DT.IVR_ATTRIBUTE2 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("H3").Value = "SELECT DT.IVR_ATTRIBUTE3 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("I3").Value = "SELECT DT.LEAD_FIELD1 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"

$ws.Range("A4").Value = "As A user I want to verify Limit Reduction"
$ws.Range("D4").Value = "Non Financial"
$ws.Range("E4").Value = "SELECT AC.DESCRIPTION FROM DC_TRANSACTION_ACTIVITY_CONFIG AC WHERE AC.TRANSACTION_TYPE_ID =(SELECT DT.TRANSACTION_TYPE_ID FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("F4").Value = "SELECT DT.IVR_ATTRIBUTE2 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("G4").Value = "SELECT DT.IVR_ATTRIBUTE2 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("H4").Value = "SELECT DT.IVR_ATTRIBUTE3 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"
$ws.Range("I4").Value = "SELECT DT.LEAD_FIELD1 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"

# Set the limit-specific columns (B = limit_type, C = new_limit) for the two
# new rows plus the updated existing row.
$ws.Range("B3").Value = "Funds Transfer to HBL Account"
$ws.Range("B2").Value = "Utility Bills and Other Payments"
$ws.Range("C3").Value = "1000"
$ws.Range("B4").Value = "Funds Transfer to Other Banks Account"
$ws.Range("C4").Value = "1"
$ws.Range("C2").Value = "23000"

# Leave the selection where it ended up after the edits.
$ws.Range("D17").Select() | Out-Null
